$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.606.54"
$ws.Range("D3").Value = "1.595.66"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.98"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.244"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.46"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("D12").Value = "1.820.34"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").Value = "1.584.47"
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.520"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.40"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "26.591.19"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "207.95"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.90"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.62%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.31"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.07%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.30"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  -1.69%  "
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.22"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("E30").Value = "  -0.96%  "
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.22"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("E33").Value = "  -2.57%  "
$ws.Range("E34").Value = "  +0.48%  "
$ws.Range("D35").Value = "1.281.82"
$ws.Range("E35").Value = "  -2.31%  "
$ws.Range("E36").Value = "  +0.92%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("E39").Value = "  +1.77%  "
$ws.Range("E41").Value = "  +1.47%  "
$ws.Range("E42").Value = "  +1.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.784"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.78"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.917"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.51%  "
$ws.Range("D46").Value = "1.732.76"
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.52"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.58"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("E50").Value = "  +3.97%  "
$ws.Range("E51").Value = "  +0.18%  "
